# filter_tags.xlsx: drop the trailing ", " left over from the tag generator
# on the two "Naval" row tag cells (C6 and B5).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C6").Value = "3D sonar technology, advanced navigation system, HISAS, Sunstone"
$ws.Range("B5").Value = "Sub-bottom profiler, Multibeam sonar, Side Scan Sonar, Single beam echo sounder"

# Update the saved cursor/selection on Sheet1 to span B24:F24 (as in the
# authored session).
$excel.Goto($ws.Range("B24:F24"))
